$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "37.819.22"
$ws.Cells.Item(2, 5).Value = "  +0.08%  "
$ws.Cells.Item(3, 4).Value = "2.082.97"
$ws.Cells.Item(3, 5).Value = "  -0.04%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "233.72"
$ws.Cells.Item(5, 5).Value = "  -0.21%  "
$ws.Cells.Item(6, 5).Value = "  +0.10%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "58.70"
$ws.Cells.Item(7, 5).Value = "  -0.52%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  +0.79%  "
$ws.Cells.Item(10, 5).Value = "  +0.00%  "
$ws.Cells.Item(11, 5).Value = "  +3.63%  "
$ws.Cells.Item(12, 5).Value = "  +1.99%  "
$ws.Cells.Item(13, 4).Value = "2.390.51"
$ws.Cells.Item(13, 5).Value = "  +0.01%  "
$ws.Cells.Item(14, 5).Value = "  +0.28%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.783"
$ws.Cells.Item(15, 5).Value = "  +1.14%  "
$ws.Cells.Item(16, 5).Value = "  +1.65%  "
$ws.Cells.Item(17, 4).Value = "2.096.04"
$ws.Cells.Item(17, 5).Value = "  +0.70%  "
$ws.Cells.Item(18, 4).Value = "37.741.99"
$ws.Cells.Item(18, 5).Value = "  +0.11%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.15"
$ws.Cells.Item(19, 5).Value = "  -0.28%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "71.54"
$ws.Cells.Item(20, 5).Value = "  +0.25%  "
$ws.Cells.Item(21, 5).Value = "  +1.62%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "230.24"
$ws.Cells.Item(22, 5).Value = "  +0.71%  "
$ws.Cells.Item(23, 5).Value = "  -0.09%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.40"
$ws.Cells.Item(24, 5).Value = "  -0.39%  "
$ws.Cells.Item(25, 5).Value = "  +1.27%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "10.02"
$ws.Cells.Item(26, 5).Value = "  +11.04%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "172.01"
$ws.Cells.Item(27, 5).Value = "  +1.11%  "
$ws.Cells.Item(28, 5).Value = "  -1.59%  "
$ws.Cells.Item(29, 5).Value = "  -0.27%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "19.55"
$ws.Cells.Item(30, 5).Value = "  +0.07%  "
$ws.Cells.Item(31, 5).Value = "  +1.24%  "
$ws.Cells.Item(32, 5).Value = "  +0.72%  "
$ws.Cells.Item(33, 5).Value = "  +0.91%  "
$ws.Cells.Item(34, 5).Value = "  -1.00%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.47"
$ws.Cells.Item(35, 5).Value = "  -1.74%  "
$ws.Cells.Item(36, 2).Value = "WEMIXToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.82"
$ws.Cells.Item(36, 5).Value = "  -0.73%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.41"
$ws.Cells.Item(37, 5).Value = "  -1.37%  "
$ws.Cells.Item(38, 5).Value = "  +0.21%  "
$ws.Cells.Item(39, 5).Value = "  +0.86%  "
$ws.Cells.Item(40, 5).Value = "  +9.75%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "101.44"
$ws.Cells.Item(41, 5).Value = "  +2.78%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0978"
$ws.Cells.Item(42, 5).Value = "  -1.23%  "
$ws.Cells.Item(43, 5).Value = "  -0.86%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "16.92"
$ws.Cells.Item(44, 5).Value = "  +5.60%  "
$ws.Cells.Item(45, 4).Value = "1.448.67"
$ws.Cells.Item(45, 5).Value = "  -0.60%  "
$ws.Cells.Item(46, 5).Value = "  -0.57%  "
$ws.Cells.Item(47, 5).Value = "  -0.04%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "4.09"
$ws.Cells.Item(48, 5).Value = "  -5.18%  "
$ws.Cells.Item(49, 5).Value = "  -0.68%  "
$ws.Cells.Item(50, 5).Value = "  -1.26%  "
$ws.Cells.Item(51, 4).Value = "2.275.30"
$ws.Cells.Item(51, 5).Value = "  +0.01%  "
